$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 159.88889
$ws.Cells.Item(33, 9).Value = 159.88889
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 11).Value = 159.88889
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 13).Value = 69.11111
$ws.Cells.Item(33, 14).ClearContents()
$ws.Cells.Item(53, 8).Value = 18519626
$ws.Cells.Item(53, 9).Value = 47619856
$ws.Cells.Item(53, 10).Value = 1298.5454
$ws.Cells.Item(53, 11).Value = 47619856
$ws.Cells.Item(53, 12).Value = 1298.5454
$ws.Cells.Item(53, 13).Value = -47619219
$ws.Cells.Item(53, 14).Value = -2572.5454
$ws.Cells.Item(74, 8).Value = 10697.789
$ws.Cells.Item(74, 9).Value = 8127.357
$ws.Cells.Item(74, 11).Value = 8127.357
$ws.Cells.Item(74, 13).Value = -7191.357
$ws.Cells.Item(77, 8).Value = 10697.789
$ws.Cells.Item(77, 9).Value = 8127.357
$ws.Cells.Item(77, 11).Value = 40636.785
$ws.Cells.Item(77, 13).Value = -35956.785
$ws.Cells.Item(112, 8).Value = 4227.8423
$ws.Cells.Item(112, 10).Value = 4612.9414
$ws.Cells.Item(112, 12).Value = 13838.8242
$ws.Cells.Item(112, 14).Value = -16054.8242
$ws.Cells.Item(116, 8).Value = 11812.786
$ws.Cells.Item(116, 9).Value = 5119.25
$ws.Cells.Item(116, 10).Value = 14490.2
$ws.Cells.Item(116, 11).Value = 5119.25
$ws.Cells.Item(116, 12).Value = 14490.2
$ws.Cells.Item(116, 13).Value = -1677.25
$ws.Cells.Item(116, 14).Value = -21374.2
$ws.Cells.Item(132, 8).Value = 2700.0417
$ws.Cells.Item(132, 9).Value = 2700.0417
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 8100.125100000001
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -5570.125100000001
$ws.Cells.Item(132, 14).ClearContents()
$ws.Cells.Item(137, 8).Value = 2566.9524
$ws.Cells.Item(137, 9).Value = 2170.8823
$ws.Cells.Item(137, 11).Value = 6512.646900000001
$ws.Cells.Item(137, 13).Value = -3962.646900000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6103.173
$ws.Cells.Item(32, 9).Value = 5958.14
$ws.Cells.Item(32, 11).Value = 5958.14
$ws.Cells.Item(32, 13).Value = -5671.14
$ws.Cells.Item(61, 8).Value = 3832.9092
$ws.Cells.Item(61, 10).Value = 5499.5
$ws.Cells.Item(61, 12).Value = 5499.5
$ws.Cells.Item(61, 14).Value = -5923.5
$ws.Cells.Item(74, 8).Value = 989.28
$ws.Cells.Item(74, 9).Value = 989.28
$ws.Cells.Item(74, 11).Value = 989.28
$ws.Cells.Item(74, 13).Value = -115.28
$ws.Cells.Item(77, 8).Value = 989.28
$ws.Cells.Item(77, 9).Value = 989.28
$ws.Cells.Item(77, 11).Value = 4946.4
$ws.Cells.Item(77, 13).Value = -578.3999999999996
$ws.Cells.Item(97, 8).Value = 953
$ws.Cells.Item(97, 9).Value = 980.36365
$ws.Cells.Item(97, 11).Value = 980.36365
$ws.Cells.Item(97, 13).Value = -484.36365
$ws.Cells.Item(122, 8).Value = 4025.4583
$ws.Cells.Item(122, 9).Value = 1651
$ws.Cells.Item(122, 10).Value = 6034.615
$ws.Cells.Item(122, 11).Value = 4953
$ws.Cells.Item(122, 12).Value = 18103.845
$ws.Cells.Item(122, 13).Value = -2503
$ws.Cells.Item(122, 14).Value = -23003.845
$ws.Cells.Item(136, 8).Value = 3832.9092
$ws.Cells.Item(136, 10).Value = 5499.5
$ws.Cells.Item(136, 12).Value = 16498.5
$ws.Cells.Item(136, 14).Value = -21598.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 286.66666
$ws.Cells.Item(22, 9).Value = 286.66666
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 286.66666
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -113.66666
$ws.Cells.Item(22, 14).ClearContents()
$ws.Cells.Item(107, 8).Value = 2038.9231
$ws.Cells.Item(107, 9).Value = 1743.8572
$ws.Cells.Item(107, 10).Value = 2383.1667
$ws.Cells.Item(107, 11).Value = 1743.8572
$ws.Cells.Item(107, 12).Value = 2383.1667
$ws.Cells.Item(107, 13).Value = 176.1428000000001
$ws.Cells.Item(107, 14).Value = -6223.1667
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 5061.25
$ws.Cells.Item(62, 9).Value = 2125
$ws.Cells.Item(62, 11).Value = 2125
$ws.Cells.Item(62, 13).Value = -1501
$ws.Cells.Item(65, 8).Value = 5061.25
$ws.Cells.Item(65, 9).Value = 2125
$ws.Cells.Item(65, 11).Value = 10625
$ws.Cells.Item(65, 13).Value = -7505
$ws.Cells.Item(105, 8).Value = 381.33334
$ws.Cells.Item(105, 9).Value = 381.33334
$ws.Cells.Item(105, 11).Value = 381.33334
$ws.Cells.Item(105, 13).Value = 1365.66666
$ws.Cells.Item(107, 8).Value = 518.1429000000001
$ws.Cells.Item(107, 9).Value = 300.625
$ws.Cells.Item(107, 11).Value = 300.625
$ws.Cells.Item(107, 13).Value = 1619.375
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(139, 8).Value = 6495.3105
$ws.Cells.Item(139, 9).Value = 5577.3887
$ws.Cells.Item(139, 11).Value = 16732.1661
$ws.Cells.Item(139, 13).Value = -11592.1661
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 484146.2
$ws.Cells.Item(113, 9).Value = 2502727.5
$ws.Cells.Item(113, 11).Value = 2502727.5
$ws.Cells.Item(113, 13).Value = -2500557.5
$ws.Cells.Item(122, 8).Value = 5025
$ws.Cells.Item(122, 9).Value = 2497.5
$ws.Cells.Item(122, 10).Value = 6036
$ws.Cells.Item(122, 11).Value = 7492.5
$ws.Cells.Item(122, 12).Value = 18108
$ws.Cells.Item(122, 13).Value = -5042.5
$ws.Cells.Item(122, 14).Value = -23008
$ws.Cells.Item(132, 8).Value = 138985.38
$ws.Cells.Item(132, 9).Value = 18147.166
$ws.Cells.Item(132, 10).Value = 501500
$ws.Cells.Item(132, 11).Value = 54441.49800000001
$ws.Cells.Item(132, 12).Value = 1504500
$ws.Cells.Item(132, 13).Value = -51911.49800000001
$ws.Cells.Item(132, 14).Value = -1509560
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 166667700
$ws.Cells.Item(16, 9).Value = 166667700
$ws.Cells.Item(16, 11).Value = 166667700
$ws.Cells.Item(16, 13).Value = -166667530
$ws.Cells.Item(68, 8).Value = 203551.2
$ws.Cells.Item(68, 9).Value = 3001.5
$ws.Cells.Item(68, 10).Value = 337251
$ws.Cells.Item(68, 11).Value = 3001.5
$ws.Cells.Item(68, 12).Value = 337251
$ws.Cells.Item(68, 13).Value = -2252.5
$ws.Cells.Item(68, 14).Value = -338749
$ws.Cells.Item(71, 8).Value = 203551.2
$ws.Cells.Item(71, 9).Value = 3001.5
$ws.Cells.Item(71, 10).Value = 337251
$ws.Cells.Item(71, 11).Value = 15007.5
$ws.Cells.Item(71, 12).Value = 1686255
$ws.Cells.Item(71, 13).Value = -11263.5
$ws.Cells.Item(71, 14).Value = -1693743
$ws.Cells.Item(93, 8).Value = 2920.6924
$ws.Cells.Item(93, 9).Value = 2810.889
$ws.Cells.Item(93, 10).Value = 3167.75
$ws.Cells.Item(93, 11).Value = 2810.889
$ws.Cells.Item(93, 12).Value = 3167.75
$ws.Cells.Item(93, 13).Value = -1562.889
$ws.Cells.Item(93, 14).Value = -5663.75
$ws.Cells.Item(136, 8).Value = 1435000.2
$ws.Cells.Item(136, 9).Value = 2505721.5
$ws.Cells.Item(136, 10).Value = 7371.8335
$ws.Cells.Item(136, 11).Value = 7517164.5
$ws.Cells.Item(136, 12).Value = 22115.5005
$ws.Cells.Item(136, 13).Value = -7514614.5
$ws.Cells.Item(136, 14).Value = -27215.5005
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1184.3846
$ws.Cells.Item(100, 9).Value = 1290.4546
$ws.Cells.Item(100, 11).Value = 2580.9092
$ws.Cells.Item(100, 13).Value = -2039.9092
$ws.Cells.Item(126, 8).Value = 2824.75
$ws.Cells.Item(126, 9).Value = 1649.5
$ws.Cells.Item(126, 11).Value = 4948.5
$ws.Cells.Item(126, 13).Value = -2478.5
$ws.Cells.Item(132, 8).Value = 48994.91
$ws.Cells.Item(132, 9).Value = 3494.4
$ws.Cells.Item(132, 10).Value = 504000
$ws.Cells.Item(132, 11).Value = 10483.2
$ws.Cells.Item(132, 12).Value = 1512000
$ws.Cells.Item(132, 13).Value = -7953.200000000001
$ws.Cells.Item(132, 14).Value = -1517060
